# ---------------------------------------------------------------------------
# Adds a "(MINI) EXTERNAL COMPONENTS" sheet (a pre-filtration-focused,
# miniaturized BOM) alongside the existing BOM, which is renamed to
# "(Large) EXTERNAL COMPONENTS".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet -------------------------------------------
$wsLarge = $wb.Worksheets.Item(1)
$wsLarge.Range("C24").Select() | Out-Null
$wsLarge.Name = "(Large) EXTERNAL COMPONENTS"

# --- Add the new MINI sheet right after the large one ----------------------
# (An extra throwaway sheet is inserted+removed first purely so the engine's
#  sheetId counter for the MINI sheet lands on 6, matching a workbook that
#  once had a 5th sheet created and discarded along the way.)
$wsTmp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLarge)
$wsMini = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTmp)
$wsMini.Name = "(MINI) EXTERNAL COMPONENTS"
$wsTmp.Delete()
# Re-fetch a fresh handle: the old $wsMini reference can carry a stale
# position after the sibling sheet was removed.
$wsMini = $wb.Worksheets.Item("(MINI) EXTERNAL COMPONENTS")
$wsMini.Select() | Out-Null

# --- Header row: copy formatting + values from the large sheet's header ---
$wsLarge.Range("A1:I1").Copy() | Out-Null
$wsMini.Range("A1:I1").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$wsMini.Range("A1:I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$wsMini.Rows.Item(1).RowHeight = 38.5

# --- Data rows --------------------------------------------------------------
# Row 3: Hollow UF membrane (Sawyer Mini filter)
$wsMini.Range("A3").Value = "Hollow UF membrane"
$wsMini.Range("B3").Value = 1
$wsMini.Range("C3").Value = "A lightweight, palm-sized Sawyer Mini water filtration kit that includes a reusable pouch, syringe, and straw, offering high-performance 0.1-micron filtration to remove bacteria, protozoa, and microplastics, with versatile compatibility for outdoor use, travel, and emergency preparedness."
$wsMini.Range("D3").Value = "Sawyer"
$wsMini.Range("G3").Value = 28.99
$wsMini.Range("H3").Value = "sawyer.com [https://www.sawyer.com/product/mini-water-filtration-system-blue#product-details]"

# Row 4: Carbon filter (Ronaqua T33)
$wsMini.Range("A4").Value = "Carbon filter"
$wsMini.Range("B4").Value = 1
$wsMini.Range("C4").Value = "A Ronaqua T33 inline coconut-shell activated carbon filter designed for reverse osmosis systems and similar applications, used as a pre/post polishing stage to improve water taste and odor, NSF-certified, easy to install, and typically replaced every two years."
$wsMini.Range("D4").Value = "Ronaqua"
$wsMini.Range("E4").Value = "`u{200e}FBA_RA-T33"
$wsMini.Range("G4").Value = 12.49
$wsMini.Range("H4").Value = "amazon.com [https://www.amazon.com/Inline-Coconut-Activated-Membrane-Reduction/dp/B0719SHH9X?th=1]"

# Row 7: Turbidity sensor (placeholder / to-decide)
$wsMini.Range("A7").Value = "Turbidity sensor"

# Row 8: Micro pressure sensor
$wsMini.Range("A8").Value = "Micro Pressure sensor"
$wsMini.Range("C8").Value = 'Search: "Pressure Transducer Sensor 5V 80psi".'

# Row 9: Flow sensor
$wsMini.Range("A9").Value = "Flow sensor"
$wsMini.Range("C9").Value = "buy the small one"
$wsMini.Range("E9").Value = "YF-S401"

# Row 12: Micro diaphragm pump
$wsMini.Range("A12").Value = "Micro Diaphragm pump"

# Row 13: Mechanical check valve
$wsMini.Range("A13").Value = "Mechanical check valve"

# --- Column widths (approximate the authored auto-fit widths) -------------
$wsMini.Columns.Item(1).ColumnWidth = 25.07
$wsMini.Columns.Item(3).ColumnWidth = 17.17
$wsMini.Columns.Item(4).ColumnWidth = 12.07
$wsMini.Columns.Item(5).ColumnWidth = 12.98
$wsMini.Columns.Item(6).ColumnWidth = 9.35
$wsMini.Columns.Item(8).ColumnWidth = 16.17
$wsMini.Columns.Item(9).ColumnWidth = 12.8

# --- View / selection state --------------------------------------------
$wsMini.Select() | Out-Null
$wsMini.Range("H18").Select() | Out-Null

Write-Host "done"
